$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("I2").Value = 114
$ws.Range("K2").Value = 142
$ws.Range("E3").Value = 146
$ws.Range("F3").Value = 137
$ws.Range("H3").Value = 155
$ws.Range("J3").Value = 231
$ws.Range("L3").Value = 245
$ws.Range("B6").Value = 375
$ws.Range("D6").Value = 414
$ws.Range("F6").Value = 537
$ws.Range("H6").Value = 440
$ws.Range("L6").Value = 430
$ws.Range("B7").Value = 503
$ws.Range("D7").Value = 646
$ws.Range("E7").Value = 697
$ws.Range("F7").Value = 776
$ws.Range("H7").Value = 717
$ws.Range("I7").Value = 833
$ws.Range("J7").Value = 795
$ws.Range("K7").Value = 890
$ws.Range("L7").Value = 827

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("F5").Value = 11
$ws.Range("F6").Value = 15

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("H3").Value = 3
$ws.Range("H6").Value = 7

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("I2").Value = 4
$ws.Range("F5").Value = 15
$ws.Range("F7").Value = 12
$ws.Range("H18").Value = 2
$ws.Range("F19").Value = 24
$ws.Range("L19").Value = 23
$ws.Range("F21").Value = 12
$ws.Range("F29").Value = 14
$ws.Range("B30").Value = 6
$ws.Range("E53").Value = 84
$ws.Range("J53").Value = 122
$ws.Range("L61").Value = 1
$ws.Range("D65").Value = 26
$ws.Range("K76").Value = 30
$ws.Range("L76").Value = 22
$ws.Range("F79").Value = 10
$ws.Range("H88").Value = 7
$ws.Range("B98").Value = 503
$ws.Range("D98").Value = 646
$ws.Range("E98").Value = 697
$ws.Range("F98").Value = 776
$ws.Range("H98").Value = 717
$ws.Range("I98").Value = 833
$ws.Range("J98").Value = 795
$ws.Range("K98").Value = 890
$ws.Range("L98").Value = 827

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("E3").Value = 15
$ws.Range("J3").Value = 36
$ws.Range("E7").Value = 84
$ws.Range("J7").Value = 122

$ws = $wb.Worksheets.Item("Sheffield & DePaul")
$ws.Range("F3").Value = 2
$ws.Range("F6").Value = 10

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K2").Value = 5
$ws.Range("L6").Value = 12
$ws.Range("K7").Value = 30
$ws.Range("L7").Value = 22

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("D5").Value = 25
$ws.Range("D6").Value = 26

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("F5").Value = 13
$ws.Range("F6").Value = 14

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("G4").Value = 2
$ws.Range("G5").Value = 2

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("F6").Value = 17
$ws.Range("L6").Value = 16
$ws.Range("F7").Value = 24
$ws.Range("L7").Value = 23

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("F5").Value = 8
$ws.Range("F6").Value = 12

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("G2").Value = 2
$ws.Range("G6").Value = 4

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("B5").Value = 4
$ws.Range("B6").Value = 6

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("F3").Value = 3
$ws.Range("F7").Value = 12

Write-Output "Applied CTA violent crime YTD update for 2025-12-07"